# Insert a new weekly price record as row 140 on the "Sandia" sheet.
# All rows that were previously 140..158 shift down to 141..159, and the
# new row 140 carries the new data point (fecha 45212, calidad "Primera",
# origin "Peru", unit "$/kilo (volumen en unidades)").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing rows 140.. down by one to make room for the new record.
$ws.Rows.Item(140).Insert()

# Populate the newly inserted row 140 with the new observation.
$ws.Cells.Item(140, 1).Value = 8
$ws.Cells.Item(140, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(140, 3).Value = "Coquimbo"
$ws.Cells.Item(140, 4).Value = 45212
$ws.Cells.Item(140, 5).Value = 4
$ws.Cells.Item(140, 6).Value = 100112028
$ws.Cells.Item(140, 7).Value = "Sandia"
$ws.Cells.Item(140, 8).Value = "Sin especificar"
$ws.Cells.Item(140, 9).Value = "Primera"
$ws.Cells.Item(140, 10).Value = 1600
$ws.Cells.Item(140, 11).Value = 800
$ws.Cells.Item(140, 12).Value = 900
$ws.Cells.Item(140, 13).Value = 850
$ws.Cells.Item(140, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(140, 15).Value = "Perú"
$ws.Cells.Item(140, 16).Value = 850
$ws.Cells.Item(140, 17).Value = 1
$ws.Cells.Item(140, 18).Value = "Hortaliza"
